$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "FlightFinder"

$ws2.Range("A1").Value = 2
$ws2.Range("B1").Value = "London"
$ws2.Range("C1").Value = "December"
$ws2.Range("D1").Value = 10
$ws2.Range("E1").Value = "Paris"
$ws2.Range("F1").Value = "December"
$ws2.Range("G1").Value = 30
$ws2.Range("H1").Value = "Unified Airlines"

$ws2.Columns.Item(8).AutoFit() | Out-Null
$ws2.Range("H1").Select()
